$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row before row 78, pushing existing rows 78-82 down to 79-83.
# (Excel copies the row-above's cell formatting onto the inserted row, so the
# new D78 already picks up the date number format used by the rest of
# column D.)
$ws.Rows.Item(78).Insert()

# Populate the newly inserted row 78 with the new record's data.
$ws.Range("A78").Value = 4
$ws.Range("B78").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C78").Value = "Los Lagos"
$ws.Range("D78").Value = 45008
$ws.Range("E78").Value = 10
$ws.Range("F78").Value = 100112043
$ws.Range("G78").Value = "Pepino dulce"
$ws.Range("H78").Value = "Cultivar IV Región"
$ws.Range("I78").Value = "Primera"
$ws.Range("J78").Value = 40
$ws.Range("K78").Value = 20000
$ws.Range("L78").Value = 20000
$ws.Range("M78").Value = 20000
$ws.Range("N78").Value = "$/bandeja 18 kilos"
$ws.Range("O78").Value = "Provincia de Limarí"
$ws.Range("P78").Value = 1111
$ws.Range("Q78").Value = 18
$ws.Range("R78").Value = "Hortaliza"
